# Insert 5 new rows of weekly "Durazno" (peach) price data at row 935,
# shifting the existing rows 935-987 down to 940-992.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A935:A939").EntireRow.Insert()

$rows = @(
  @{Row=935; D=44931; K="Carson";    L="Especial"; M=170; N=16800; O=16800; P=16800; Q="`$/caja 14 kilos empedrada"; R="Provincia de Los Andes";    S=1200; T=14},
  @{Row=936; D=44931; K="Carson";    L="Primera";  M=150; N=14000; O=14000; P=14000; Q="`$/caja 14 kilos empedrada"; R="Provincia de Los Andes";    S=1000; T=14},
  @{Row=937; D=44931; K="Carson";    L="Segunda";  M=180; N=11200; O=11200; P=11200; Q="`$/caja 14 kilos empedrada"; R="Provincia de Los Andes";    S=800;  T=14},
  @{Row=938; D=44931; K="Rich Lady"; L="Especial"; M=350; N=16000; O=16000; P=16000; Q="`$/caja 16 kilos granel";    R="Región de O'Higgins";       S=1000; T=16},
  @{Row=939; D=44931; K="Rich Lady"; L="Primera";  M=280; N=12800; O=12800; P=12800; Q="`$/caja 16 kilos granel";    R="Región de O'Higgins";       S=800;  T=16}
)

foreach ($r in $rows) {
  $rowNum = $r.Row
  $ws.Range("A$rowNum").Value = 9
  $ws.Range("B$rowNum").Value = "Vega Central Mapocho de Santiago"
  $ws.Range("C$rowNum").Value = "Metropolitana"
  $ws.Range("D$rowNum").Value = $r.D
  $ws.Range("E$rowNum").Value = 13
  $ws.Range("F$rowNum").Value = "Fruta"
  $ws.Range("G$rowNum").Value = 100103
  $ws.Range("H$rowNum").Value = "Frutos de hueso (carozo)"
  $ws.Range("I$rowNum").Value = 100103004
  $ws.Range("J$rowNum").Value = "Durazno"
  $ws.Range("K$rowNum").Value = $r.K
  $ws.Range("L$rowNum").Value = $r.L
  $ws.Range("M$rowNum").Value = $r.M
  $ws.Range("N$rowNum").Value = $r.N
  $ws.Range("O$rowNum").Value = $r.O
  $ws.Range("P$rowNum").Value = $r.P
  $ws.Range("Q$rowNum").Value = $r.Q
  $ws.Range("R$rowNum").Value = $r.R
  $ws.Range("S$rowNum").Value = $r.S
  $ws.Range("T$rowNum").Value = $r.T
}
